$wb = $excel.ActiveWorkbook

# --- Sheet "Eetu Pihamäki": add the missing work-log entry on row 23 ---
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# Date: 29.10.2018
$ws.Range("A23").Value = 43402
# Start time 17:15, end time 19:20
$ws.Range("B23").Value = 0.71875
$ws.Range("C23").Value = 0.80555555555555547
# Sprint number
$ws.Range("E23").Value = 3
# Task description
$ws.Range("F23").Value = "1 h Projektipäällikön hallinnollisia tehtäviä (Pöytäkirja,  tuntikirjanpito jne.) 50 min Windows Server 2016 aktivointiavain toimimaan --> Microsoft Supportille soitto ja chat.  https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2029.10.2018.txt"

# Row grows to fit the wrapped task text, same as the neighbouring rows
$ws.Rows.Item(23).RowHeight = 90

# Scroll the sheet view down one row (topLeftCell A19 -> A20)
$win = $excel.ActiveWindow
$win.ScrollRow = 20
